$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 08:16"
$ws.Range("B15").Value = 4477
$ws.Range("C15").Value = 3
$ws.Range("E15").Value = 4447
$ws.Range("F15").Value = 16
$ws.Range("B21").Value = 2059
$ws.Range("C21").Value = 13
$ws.Range("E21").Value = 2010
$ws.Range("F21").Value = 110
$ws.Range("G21").Value = 6
$ws.Range("H21").Value = 33
$ws.Range("E24").Value = 1344
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 15
$ws.Range("B28").Value = 1140
$ws.Range("C28").Value = 12
$ws.Range("D28").Value = 285
$ws.Range("E28").Value = 813
$ws.Range("F28").Value = 54
$ws.Range("B67").Value = 202
$ws.Range("C67").Value = 1
$ws.Range("E67").Value = 196
$ws.Range("A79").Value = "Bosnia y Herzegovina"
$ws.Range("B79").Value = 137
$ws.Range("C79").Value = 1
$ws.Range("D79").Value = 2
$ws.Range("E79").Value = 134
$ws.Range("H79").Value = 1
$ws.Range("A80").Value = "Republica de Macedonia"
$ws.Range("D80").Value = 1
$ws.Range("H80").Value = 2
$ws.Range("A110").Value = "Kirguistan"
$ws.Range("C110").Value = 26
$ws.Range("D110").Value = 0
$ws.Range("E110").Value = 42
$ws.Range("H110").Value = 0
$ws.Range("A111").Value = "Afganistan"
$ws.Range("B111").Value = 42
$ws.Range("C111").Value = 2
$ws.Range("D111").Value = 1
$ws.Range("E111").Value = 40
$ws.Range("F111").Value = 0
$ws.Range("A112").Value = "Cuba"
$ws.Range("D112").Value = 0
$ws.Range("E112").Value = 39
$ws.Range("F112").Value = 3
$ws.Range("A113").Value = "Nigeria"
$ws.Range("B113").Value = 40
$ws.Range("D113").Value = 2
$ws.Range("E113").Value = 37
$ws.Range("H113").Value = 1
$ws.Range("A114").Value = "Ruanda"
$ws.Range("E114").Value = 36
$ws.Range("H114").Value = 0
$ws.Range("A116").Value = "Consejo Danes para los Refugiados"
$ws.Range("B116").Value = 36
$ws.Range("C116").Value = 0
$ws.Range("E116").Value = 34
$ws.Range("H116").Value = 2
$ws.Range("A120").Value = "Guam"
$ws.Range("B120").Value = 29
$ws.Range("C120").Value = 0
$ws.Range("H120").Value = 1
$ws.Range("A121").Value = "Bolivia"
$ws.Range("B121").Value = 28
$ws.Range("C121").Value = 1
$ws.Range("E121").Value = 28
$ws.Range("H121").Value = 0
$ws.Range("A122").Value = "Montenegro"
$ws.Range("C122").Value = 0
$ws.Range("E122").Value = 26
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 1
$ws.Range("A123").Value = "Paraguay"
$ws.Range("C123").Value = 5
$ws.Range("F123").Value = 1
$ws.Range("G123").Value = 1
$ws.Range("A124").Value = "Ghana"
$ws.Range("B124").Value = 27
$ws.Range("D124").Value = 0
$ws.Range("E124").Value = 25
$ws.Range("H124").Value = 2
$ws.Range("A125").Value = "Costa de Marfil"
$ws.Range("D125").Value = 2
$ws.Range("E125").Value = 23
$ws.Range("A126").Value = "Macao"
$ws.Range("B126").Value = 25
$ws.Range("D126").Value = 10
$ws.Range("E126").Value = 15
$ws.Range("A127").Value = "Mayotte"
$ws.Range("B127").Value = 24
$ws.Range("C127").Value = 0
$ws.Range("E127").Value = 24
$ws.Range("A128").Value = "Polinesia Francesa"
$ws.Range("C128").Value = 5
$ws.Range("D128").Value = 0
$ws.Range("E128").Value = 23
$ws.Range("A129").Value = "Monaco"
$ws.Range("B129").Value = 23
$ws.Range("D129").Value = 1
$ws.Range("E129").Value = 22
$ws.Range("H129").Value = 0
$ws.Range("A130").Value = "Guyana"
$ws.Range("A131").Value = "Guatemala"
$ws.Range("D131").Value = 0
$ws.Range("E131").Value = 19
$ws.Range("H131").Value = 1
$ws.Range("A132").Value = "Guayana Francesa"
$ws.Range("B132").Value = 20
$ws.Range("D132").Value = 6
$ws.Range("E132").Value = 14
$ws.Range("H132").Value = 0
$ws.Range("A133").Value = "Jamaica"
$ws.Range("B133").Value = 19
$ws.Range("D133").Value = 2
$ws.Range("E133").Value = 16
$ws.Range("H133").Value = 1
$ws.Range("A134").Value = "Togo"
$ws.Range("B134").Value = 18
$ws.Range("E134").Value = 18
$ws.Range("A135").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("C135").Value = 0
$ws.Range("A136").Value = "Madagascar"
$ws.Range("C136").Value = 5
$ws.Range("A137").Value = "Barbados"
$ws.Range("B137").Value = 17
$ws.Range("E137").Value = 17
$ws.Range("A146").Value = "Guinea Ecuatorial"
$ws.Range("A147").Value = "Uganda"
$ws.Range("A148").Value = "San Martin (Parte Francesa)"
$ws.Range("A149").Value = "Nueva Caledonia"
$ws.Range("A152").Value = "Bermudas"
$ws.Range("A154").Value = "Haiti"
$ws.Range("A159").Value = "Suazilandia"
$ws.Range("A161").Value = "Guinea"
$ws.Range("A162").Value = "Namibia"
$ws.Range("A163").Value = "Bahamas"
$ws.Range("A164").Value = "Congo"
$ws.Range("A178").Value = "Birmania"
$ws.Range("A181").Value = "Nicaragua"
$ws.Range("A183").Value = "Butan"
